$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New firm-usage rows for Firm Name / Duration / Lawyers (A:C), rows 2-52.
$rows = @(
    @("Spencer West", "27sec", "1"),
    @("Conyers", "43sec", "1"),
    @("Krogerus", "11sec", "1"),
    @("Poulschmith", "29sec", "1"),
    @("Gianni And Origoni", "13sec", "1"),
    @("Pearl Cohen", "12sec", "1"),
    @("Consortium Legal", "8sec", "1"),
    @("GÖRG", "18sec", "1"),
    @("Beauchamps", "6sec", "1"),
    @("Keystone Law", "12sec", "1"),
    @("Latham And Watkins", "20sec", "1"),
    @("Howse Williams", "17sec", "1"),
    @("HFW", "20sec", "1"),
    @("Deacons", "13sec", "1"),
    @("Schoenherr", "2min 50sec", "1"),
    @("Pedersoli", "7sec", "1"),
    @("Gornitzky And Co", "12sec", "1"),
    @("Walkers", "42sec", "1"),
    @("Hannes Snellman", "17sec", "1"),
    @("BonelliErede", "21sec", "1"),
    @("Latam Lex", "8sec", "1"),
    @("White and Case", "14sec", "1"),
    @("Hill Dickinson", "8sec", "1"),
    @("Lex Caribbean", "14sec", "1"),
    @("Carey Olsen", "19sec", "1"),
    @("Aron Tadmor Levy", "19sec", "1"),
    @("Myers Fletcher And Gordon", "16sec", "1"),
    @("Fangda Partners", "19sec", "1"),
    @("Matheson", "11sec", "1"),
    @("Gide Loyrette Nouel", "11sec", "1"),
    @("Remfry And Sagar", "21sec", "1"),
    @("Asafo And Co", "16sec", "1"),
    @("Havel Partners", "11sec", "1"),
    @("Kinstellar", "29sec", "1"),
    @("Pulegal", "11sec", "1"),
    @("Harneys", "2min 3sec", "1"),
    @("Simmons And Simmons", "51sec", "1"),
    @("Wolf Theiss", "48sec", "1"),
    @("Longan Law", "10sec", "1"),
    @("Samvad Partners", "4sec", "1"),
    @("Winston And Strawn", "8sec", "1"),
    @("Tanner DeWitt", "6min 39sec", "1"),
    @("Al Tamimi", "23sec", "1"),
    @("DGKV", "15sec", "1"),
    @("Kromann Reumert", "15sec", "1"),
    @("Dittmar And Indrenius", "9sec", "1"),
    @("Pavia And Ansaldo", "14sec", "1"),
    @("DSK Legal", "11sec", "1"),
    @("Carey Olsen", "17sec", "1"),
    @("EBN", "14sec", "1"),
    @("Borenius", "29sec", "1")
)

# A2:C52 in the source sheet is stored as text (numberStoredAsText); switch
# to a text number format before assignment so a purely-numeric string like
# "1" (the Lawyers column) is not auto-coerced into a numeric cell by COM,
# then restore "General" so the cell style matches the source workbook.
$range = $ws.Range("A2:C52")
$range.NumberFormat = "@"

$r = 2
foreach ($item in $rows) {
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $r = $r + 1
}

$range.NumberFormat = "General"
